# Adds season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the style of the last existing header cell (AC1) onto the
# three new header cells, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-42: every row gets the same season record (75-87-0).
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 75
    $ws.Cells.Item($row, 31).Value = 87
    $ws.Cells.Item($row, 32).Value = 0
}
